$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns I (I0) and J (IF)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header formatting (bold, bordered, centered) used by the other
# header cells (e.g. H1) by copying its format onto the new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# New data values for columns I and J (rows 2-10)
$dataI = @(2, 3, 3, 1, 7, 4, 1, 4, 3)
$dataJ = @(4, 4, 5, 3, 9, 6, 5, 6, 4)

for ($i = 0; $i -lt 9; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $dataI[$i]
    $ws.Cells.Item($row, 10).Value = $dataJ[$i]
}
